$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header of column A from "Gen" to "MaxFES"
$ws.Range("A1").Value = "MaxFES"

# 2. Update the "MaxFES" values in column A (row 2 stays 0, rows 3-14 change)
$ws.Range("A3").Value  = 0.001
$ws.Range("A4").Value  = 0.01
$ws.Range("A5").Value  = 0.1
$ws.Range("A6").Value  = 0.2
$ws.Range("A7").Value  = 0.3
$ws.Range("A8").Value  = 0.4
$ws.Range("A9").Value  = 0.5
$ws.Range("A10").Value = 0.6
$ws.Range("A11").Value = 0.7
$ws.Range("A12").Value = 0.8
$ws.Range("A13").Value = 0.9
$ws.Range("A14").Value = 1

# 3. Delete the "Run 50" column (column AZ, the 52nd column).
#    This shifts the old "Mean" column (BA) left into AZ, and Excel
#    automatically drops the now-unused "Run 50" shared string.
$ws.Columns.Item(52).Delete()

# 4. The shifted-in "Mean" column (now AZ) still holds the old cached
#    mean (computed over all 51 runs incl. Run 50). Recompute it as the
#    mean over the remaining 50 runs (columns B:AY) for every data row.
$ws.Range("AZ2").Value  = 71.58003983
$ws.Range("AZ3").Value  = 48.85323682
$ws.Range("AZ4").Value  = 4.34840451
$ws.Range("AZ5").Value  = 0.24762849
$ws.Range("AZ6").Value  = 0.20722963
$ws.Range("AZ7").Value  = 0.18589199
$ws.Range("AZ8").Value  = 0.17067023
$ws.Range("AZ9").Value  = 0.1617451
$ws.Range("AZ10").Value = 0.15468213
$ws.Range("AZ11").Value = 0.14678065
$ws.Range("AZ12").Value = 0.14069703
$ws.Range("AZ13").Value = 0.13787067
$ws.Range("AZ14").Value = 0.13415723
